$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40..59 down to 41..60
$ws.Rows(40).Insert()

# Populate the newly inserted row 40 with the new record's data
$ws.Cells.Item(40, 1).Value = 9
$ws.Cells.Item(40, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(40, 3).Value = "Metropolitana"
$ws.Cells.Item(40, 4).Value = 44452
$ws.Cells.Item(40, 5).Value = 13
$ws.Cells.Item(40, 6).Value = 100112022
$ws.Cells.Item(40, 7).Value = "Arveja Verde"
$ws.Cells.Item(40, 8).Value = "Perfection"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 25
$ws.Cells.Item(40, 11).Value = 37000
$ws.Cells.Item(40, 12).Value = 38000
$ws.Cells.Item(40, 13).Value = 37480
$ws.Cells.Item(40, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(40, 16).Value = 1499
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"
